$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "303.92"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "6.36%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.94"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "8.62%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.263"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "3.87%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07501"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "11.76%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.833"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "6.97%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.765"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "9.44%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.485"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "6.34%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9143"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.41%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01660"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2,469.05%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1689"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "6.99%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07526"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "5.52%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08068"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "5.56%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.02991"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "2.74%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09906"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "10.30%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001493"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-5.82%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04560"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.73%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006408"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "3.39%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.33%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.03%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.39%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.1344"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2.33%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.487"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "13.93%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1630"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "4.55%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001217"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.17%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004445"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "1.76%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001398"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "19.52%"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001739"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "7.48%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04500"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "6.37%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007212"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "5.92%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1348"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002247"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.77%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01301"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "3.28%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006226"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "8.18%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.869"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-4.78%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.01298"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-13.49%"
